$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.785.14"
$ws.Range("E2").Value = "  +2.24%  "

$ws.Range("D3").Value = "'1.656.79"
$ws.Range("E3").Value = "  +3.13%  "

$ws.Range("D4").Value = "'0.9966"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "'306.83"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").Value = "'0.9974"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").Value = "'0.3776"
$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("D8").Value = "'52.73"
$ws.Range("E8").Value = "  +0.61%  "

$ws.Range("D9").Value = "'0.3683"
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("D10").Value = "'1.278"
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("D11").Value = "'0.08184"
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "'0.9968"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "'23.27"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").Value = "'6.744"
$ws.Range("E14").Value = "  +2.19%  "

$ws.Range("D15").Value = "'0.00001280"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("D16").Value = "'7.426"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "'1.660.43"
$ws.Range("E17").Value = "  +3.32%  "

$ws.Range("D18").Value = "'95.54"
$ws.Range("E18").Value = "  +1.67%  "

$ws.Range("D19").Value = "'0.06918"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "'18.55"
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("D21").Value = "'6.622"
$ws.Range("E21").Value = "  +1.21%  "

$ws.Range("D22").Value = "'0.9981"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").Value = "'23.773.37"
$ws.Range("E23").Value = "  +2.24%  "

$ws.Range("D24").Value = "'13.03"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").Value = "'3.249"
$ws.Range("E25").Value = "  +5.78%  "

$ws.Range("D26").Value = "'2.434"
$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("D27").Value = "'21.51"
$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("D28").Value = "'151.38"

$ws.Range("D29").Value = "'5.340"
$ws.Range("E29").Value = "  +1.01%  "

$ws.Range("D30").Value = "'137.58"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("D31").Value = "'2.330"
$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("D32").Value = "'1.846.08"
$ws.Range("E32").Value = "  +3.64%  "

$ws.Range("D33").Value = "'6.940"
$ws.Range("E33").Value = "  +2.59%  "

$ws.Range("D34").Value = "'11.14"
$ws.Range("E34").Value = "  +7.30%  "

$ws.Range("D35").Value = "'0.9803"
$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("D36").Value = "'0.02889"
$ws.Range("E36").Value = "  +4.57%  "

$ws.Range("D37").Value = "'6.438"
$ws.Range("E37").Value = "  +5.21%  "

$ws.Range("D38").Value = "'0.2606"
$ws.Range("E38").Value = "  +3.31%  "

$ws.Range("D39").Value = "'0.07381"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").Value = "'0.08913"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("D41").Value = "'0.7268"
$ws.Range("E41").Value = "  +2.21%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.386"
$ws.Range("E42").Value = "  -3.02%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'12.84"
$ws.Range("E43").Value = "  +2.55%  "

$ws.Range("E44").Value = "  +5.31%  "

$ws.Range("D45").Value = "'0.6681"
$ws.Range("E45").Value = "  +2.03%  "

$ws.Range("D46").Value = "'2.400"
$ws.Range("E46").Value = "  +2.54%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'4.034"
$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'0.9970"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("D49").Value = "'0.08073"
$ws.Range("E49").Value = "  +1.65%  "

$ws.Range("D50").Value = "'1.236"
$ws.Range("E50").Value = "  +2.04%  "

$ws.Range("D51").Value = "'129.01"
$ws.Range("E51").Value = "  -3.33%  "
